$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 61, shifting the existing
# rows 61..90 down to 62..91 (this matches the dimension change
# A1:T90 -> A1:T91 and the observed row-data shift in the diff).
$ws.Rows("61").Insert()

# Populate the newly inserted row 61 with the new weekly price entry.
$ws.Cells.Item(61, 1).Value = 11
$ws.Cells.Item(61, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(61, 3).Value = "Bíobío"
$ws.Cells.Item(61, 4).Value = 44510
$ws.Cells.Item(61, 5).Value = 8
$ws.Cells.Item(61, 6).Value = "Fruta"
$ws.Cells.Item(61, 7).Value = 100102
$ws.Cells.Item(61, 8).Value = "Cítricos"
$ws.Cells.Item(61, 9).Value = 100102004
$ws.Cells.Item(61, 10).Value = "Mandarina"
$ws.Cells.Item(61, 11).Value = "Murcott"
$ws.Cells.Item(61, 12).Value = "Primera"
$ws.Cells.Item(61, 13).Value = 300
$ws.Cells.Item(61, 14).Value = 7000
$ws.Cells.Item(61, 15).Value = 7500
$ws.Cells.Item(61, 16).Value = 7250
$ws.Cells.Item(61, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(61, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(61, 19).Value = 403
$ws.Cells.Item(61, 20).Value = 18
